{"js": "// Replace the two-digit multiplication problems/answers throughout the\n// document body (including any tables) with the updated values from the\n// commit. Each source string is unique in the document, so a targeted\n// search-and-replace per pair is safe and preserves all run formatting\n// (font, size, etc.) on the matched run.\nconst replacements = [\n  [\"67\u00d765=4355\", \"81\u00d718=1458\"],\n  [\"25\u00d724=600\", \"54\u00d769=3726\"],\n  [\"93\u00d712=1116\", \"30\u00d790=2700\"],\n  [\"42\u00d789=3738\", \"56\u00d792=5152\"],\n  [\"91\u00d722=2002\", \"99\u00d774=7326\"],\n  [\"73\u00d724=1752\", \"82\u00d797=7954\"],\n  [\"95\u00d729=2755\", \"62\u00d765=4030\"],\n  [\"35\u00d727=945\", \"68\u00d718=1224\"],\n  [\"39\u00d719=741\", \"98\u00d725=2450\"],\n  [\"90\u00d762=5580\", \"94\u00d780=7520\"],\n  [\"16\u00d736=576\", \"92\u00d766=6072\"],\n  [\"33\u00d753=1749\", \"23\u00d762=1426\"],\n  [\"59\u00d782=4838\", \"74\u00d735=2590\"],\n  [\"80\u00d716=1280\", \"21\u00d712=252\"],\n  [\"38\u00d743=1634\", \"60\u00d736=2160\"],\n  [\"29\u00d794=2726\", \"42\u00d768=2856\"],\n  [\"86\u00d771=6106\", \"45\u00d729=1305\"],\n  [\"13\u00d739=507\", \"78\u00d753=4134\"],\n  [\"93\u00d786=7998\", \"65\u00d767=4355\"],\n  [\"26\u00d743=1118\", \"22\u00d714=308\"],\n  [\"64\u00d798=6272\", \"92\u00d773=6716\"],\n  [\"43\u00d732=1376\", \"70\u00d794=6580\"],\n  [\"64\u00d725=1600\", \"87\u00d742=3654\"],\n  [\"77\u00d777=5929\", \"68\u00d751=3468\"],\n  [\"18\u00d742=756\", \"76\u00d738=2888\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems/answers throughout the\n# document body (including any tables) with the updated values from the\n# commit. Each source string is unique in the document, so a targeted\n# Find/Replace per pair is safe and preserves all run formatting\n# (font, size, etc.) on the matched run.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"67\u00d765=4355\", \"81\u00d718=1458\"),\n  @(\"25\u00d724=600\", \"54\u00d769=3726\"),\n  @(\"93\u00d712=1116\", \"30\u00d790=2700\"),\n  @(\"42\u00d789=3738\", \"56\u00d792=5152\"),\n  @(\"91\u00d722=2002\", \"99\u00d774=7326\"),\n  @(\"73\u00d724=1752\", \"82\u00d797=7954\"),\n  @(\"95\u00d729=2755\", \"62\u00d765=4030\"),\n  @(\"35\u00d727=945\", \"68\u00d718=1224\"),\n  @(\"39\u00d719=741\", \"98\u00d725=2450\"),\n  @(\"90\u00d762=5580\", \"94\u00d780=7520\"),\n  @(\"16\u00d736=576\", \"92\u00d766=6072\"),\n  @(\"33\u00d753=1749\", \"23\u00d762=1426\"),\n  @(\"59\u00d782=4838\", \"74\u00d735=2590\"),\n  @(\"80\u00d716=1280\", \"21\u00d712=252\"),\n  @(\"38\u00d743=1634\", \"60\u00d736=2160\"),\n  @(\"29\u00d794=2726\", \"42\u00d768=2856\"),\n  @(\"86\u00d771=6106\", \"45\u00d729=1305\"),\n  @(\"13\u00d739=507\", \"78\u00d753=4134\"),\n  @(\"93\u00d786=7998\", \"65\u00d767=4355\"),\n  @(\"26\u00d743=1118\", \"22\u00d714=308\"),\n  @(\"64\u00d798=6272\", \"92\u00d773=6716\"),\n  @(\"43\u00d732=1376\", \"70\u00d794=6580\"),\n  @(\"64\u00d725=1600\", \"87\u00d742=3654\"),\n  @(\"77\u00d777=5929\", \"68\u00d751=3468\"),\n  @(\"18\u00d742=756\", \"76\u00d738=2888\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  # wdFindContinue=1, wdReplaceAll=2\n  $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
